# Add a new column R (year 2021) to the malaria incidence table on sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new year column
$ws.Range("R3").Value = 2021

# Data column for 2021. Most rows keep the "-" placeholder used throughout
# the sheet for missing data; a couple of rows get real figures.
$ws.Range("R4").Value = 0.00029886145739191973
$ws.Range("R5").Value = "-"
$ws.Range("R6").Value = "-"
$ws.Range("R7").Value = "-"
$ws.Range("R8").Value = "-"
$ws.Range("R9").Value = "-"
$ws.Range("R10").Value = "-"
$ws.Range("R11").Value = "-"
$ws.Range("R12").Value = 0.0018411781330637848
$ws.Range("R13").Value = "-"

# Match formatting of the neighbouring Q column cells for the new column.
$ws.Range("R3").Style = $ws.Range("Q3").Style
$ws.Range("R4").Style = $ws.Range("N4").Style
$ws.Range("R5:R11").Style = $ws.Range("Q5").Style
$ws.Range("R12").NumberFormat = "0.0"
$ws.Range("R13").Style = $ws.Range("Q13").Style

# Update the selection to mirror the saved file's cursor position.
$ws.Range("S4").Select()
